$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new summary rows below the data table -------------------

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Style the four summary values: bold, 12pt, vertically centered -------
# Build the style once on B14 and fan it out with copy/paste so only a
# single new font + cellXf combination is introduced (matches how Excel
# itself collapses identical formatting into one style record).
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Taller rows for the new summary section
$ws.Range("A14:B17").RowHeight = 15.6

# --- Page setup tweaks ------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore the selection that was active when the file was last saved ---
$ws.Range("C16").Select()
